# Apply "First Game results - official" update to the Milano season2 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Andrea Limonta -> Luca Stoppi
$ws.Range("A2").Value = "Luca Stoppi"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 18
$ws.Range("O2").Value = 18

# Row 3: Cerro -> Andrea Limonta
$ws.Range("A3").Value = "Andrea Limonta"
$ws.Range("B3").Value = 2
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1

# Row 4: Luca Stoppi -> Giovanni Aiello
$ws.Range("A4").Value = "Giovanni Aiello"
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 14
$ws.Range("O4").Value = 14

# Row 5: Giovanni Aiello -> Maurizio
$ws.Range("A5").Value = "Maurizio"
$ws.Range("B5").Value = 3

# Row 6: Maurizio -> Cerro
$ws.Range("A6").Value = "Cerro"
$ws.Range("B6").Value = 5
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 13

# Row 7: Mazzu rank update
$ws.Range("B7").Value = 5
